# Swap the order of names in the "Recorded By" column (G) wherever the
# session was recorded by both the System and dnasr281@gmail.com, i.e.
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$first = $null
$found = $ws.Cells.Find($target)
while ($found -ne $null) {
    if ($first -eq $null) {
        $first = $found.Address()
    } elseif ($found.Address() -eq $first) {
        break
    }
    $found.Value = $replacement
    $found = $ws.Cells.FindNext($found)
}
